# Apply the "property boat&car done" edit to the 汽車 (car) worksheet.
# The car sheet previously had its header row (row 1) duplicating the
# first data row, was missing a "capacity" (排氣量) column, and was
# missing the metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that the other
# property sheets (land/building/stock/...) already had.
#
# This script fixes the header row, inserts the new "capacity" column
# (C), and fills in the missing metadata columns (H:N) for the existing
# data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Row 1: proper column headers ----
$ws.Cells.Item(1, 2).Value2 = "name"
$ws.Cells.Item(1, 3).Value2 = "capacity"
$ws.Cells.Item(1, 4).Value2 = "owner"
$ws.Cells.Item(1, 5).Value2 = "register_date"
$ws.Cells.Item(1, 6).Value2 = "register_reason"
$ws.Cells.Item(1, 7).Value2 = "acquire_value"
$ws.Cells.Item(1, 8).Value2 = "property_category"
$ws.Cells.Item(1, 9).Value2 = "category"
$ws.Cells.Item(1, 10).Value2 = "date"
$ws.Cells.Item(1, 11).Value2 = "legislator_name"
$ws.Cells.Item(1, 12).Value2 = "legislator_id"
$ws.Cells.Item(1, 13).Value2 = "source_file"
$ws.Cells.Item(1, 14).Value2 = "index"

# Give the newly added header cells (H1:N1) the same look (bold, centered,
# bordered) as the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# ---- Row 2: fill in the newly added metadata columns ----
$ws.Cells.Item(2, 8).Value2 = "land"
$ws.Cells.Item(2, 9).Value2 = "normal"
# Force the "date" column to stay plain text (e.g. "2012-04-25") instead
# of being auto-converted into a date serial number.
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value2 = "2012-04-25"
$ws.Cells.Item(2, 11).Value2 = "陳碧涵"
$ws.Cells.Item(2, 12).Value2 = 1752
$ws.Cells.Item(2, 13).Value2 = "tmpd4df1"
$ws.Cells.Item(2, 14).Value2 = 44

# Match the plain data-row styling used by the rest of row 2 (this also
# resets the number format that was forced above back to the shared
# "normal" data style).
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
